# Weekly fruit/hortaliza update: insert a new daily price record as the
# most recent entry (row 329) for "Femacal de La Calera" - Acelga, pushing
# the previously existing rows 329-347 down to 330-348.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row before row 329, shifting rows 329:347 -> 330:348.
$ws.Rows.Item(329).Insert()

# Populate the newly inserted row 329 with the new week's record.
$ws.Cells.Item(329, 1).Value = 3
$ws.Cells.Item(329, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(329, 3).Value = "Coquimbo"
$ws.Cells.Item(329, 4).Value = 44706
$ws.Cells.Item(329, 5).Value = 5
$ws.Cells.Item(329, 6).Value = 100112009
$ws.Cells.Item(329, 7).Value = "Acelga"
$ws.Cells.Item(329, 8).Value = "Sin especificar"
$ws.Cells.Item(329, 9).Value = "Primera"
$ws.Cells.Item(329, 10).Value = 280
$ws.Cells.Item(329, 11).Value = 3000
$ws.Cells.Item(329, 12).Value = 3300
$ws.Cells.Item(329, 13).Value = 3171
$ws.Cells.Item(329, 14).Value = "$/docena de atados (6 kilos)"
$ws.Cells.Item(329, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(329, 16).Value = 528
$ws.Cells.Item(329, 17).Value = 6
$ws.Cells.Item(329, 18).Value = "Hortaliza"
